$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update runMode column for rows 4-7 from "y" to "n"
$ws.Range("C4").Value = "n"
$ws.Range("C5").Value = "n"
$ws.Range("C6").Value = "n"
$ws.Range("C7").Value = "n"

# Move the active selection to C8
$ws.Range("C8").Select()
